$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.442.67"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.572.40"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'1.003"

$ws.Range("D6").Value = "'291.13"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "'0.3734"
$ws.Range("E7").Value = "  -0.91%  "

$ws.Range("D8").Value = "'49.88"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.3393"

$ws.Range("D10").Value = "'0.07559"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("D11").Value = "'1.134"
$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "'21.32"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "'5.984"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").Value = "'6.943"
$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").Value = "1.592.68"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "'0.00001119"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "'90.86"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "'0.06727"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").Value = "'6.277"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "'16.33"
$ws.Range("E22").Value = "  -3.16%  "

$ws.Range("D23").Value = "'12.12"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").Value = "22.460.66"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").Value = "'2.331"
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("D26").Value = "'2.632"
$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").Value = "'20.13"
$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("D28").Value = "'148.11"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").Value = "'5.021"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").Value = "'125.44"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("D31").Value = "1.764.67"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("D32").Value = "'1.045"
$ws.Range("E32").Value = "  +5.46%  "

$ws.Range("D33").Value = "'6.117"
$ws.Range("E33").Value = "  -1.05%  "

$ws.Range("D35").Value = "'9.738"
$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").Value = "'0.08376"
$ws.Range("E36").Value = "  -2.62%  "

$ws.Range("D37").Value = "'1.384"
$ws.Range("E37").Value = "  +4.73%  "

$ws.Range("E38").Value = "  -3.24%  "

$ws.Range("D39").Value = "'0.2284"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("D40").Value = "'0.06512"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("D41").Value = "'5.443"
$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("D42").Value = "'11.21"
$ws.Range("E42").Value = "  -2.71%  "

$ws.Range("D43").Value = "'0.6203"
$ws.Range("E43").Value = "  -3.70%  "

$ws.Range("D45").Value = "'13.89"
$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").Value = "'3.811"

$ws.Range("D47").Value = "'0.5775"
$ws.Range("E47").Value = "  -3.91%  "

$ws.Range("D48").Value = "'129.72"
$ws.Range("E48").Value = "  +3.22%  "

$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").Value = "'1.213"
$ws.Range("E50").Value = "  -6.87%  "

$ws.Range("D51").Value = "'0.07307"
$ws.Range("E51").Value = "  -0.30%  "
